$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 4 (year 2030) with new probability values
$ws.Range("B4").Value = 0.62
$ws.Range("C4").Value = 0.3
$ws.Range("D4").Value = 0.08

# Add new row 5 (year 2040) with the same probability values as 2030
$ws.Range("A5").Value = 2040
$ws.Range("B5").Value = 0.62
$ws.Range("C5").Value = 0.3
$ws.Range("D5").Value = 0.08

# Add new row 6 (year 2050) re-using the original 2035 probability values
# (copy straight from row 3, which still holds the untouched 0.81/0.15/0.04
# values, so the stored floating point bit patterns match exactly)
$ws.Range("A6").Value = 2050
$ws.Range("B6").Value = $ws.Range("B3").Value2
$ws.Range("C6").Value = $ws.Range("C3").Value2
$ws.Range("D6").Value = $ws.Range("D3").Value2
